$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number. Every data row
# (2..332) had it bumped from 45206 (2023-10-07) to 45208 (2023-10-09).
$lastCell = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162)
$lastRow = $lastCell.Row

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45208
}
